# Updates the loading_percent results for Case_1_84 (380 kV case) on Sheet1.
# Rows 2-25 correspond to the 24 time steps; columns B,D:L,N:O hold the per-line
# loading percentages that were recomputed for the new case (columns A, C, M are
# unchanged index / zero columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Key = 1-based column number, Value = array of new values for rows 2..25
$data = @{
    2 = @(13.768360054235; 13.63173930748824; 13.54917428245344; 13.51589344958914; 13.51039019477863; 13.54872392415938; 13.72099908583542; 14.06783564962402; 14.32615222451117; 14.44400111051704; 14.48864085488417; 14.47902687832179; 14.4476735938362; 14.42846942878078; 14.31845427516861; 14.25102480668211; 14.21227522770187; 14.19916223740831; 14.25819952410296; 14.45688274616652; 14.58679280289332; 14.51746367630768; 14.25495578015786; 13.97324841812893) # column B
    4 = @(10.30614043403994; 10.31189964328614; 10.31677790713182; 10.31910417764741; 10.31951091162707; 10.31680790901797; 10.30784819506214; 10.30088722790366; 10.30218393462986; 10.30415231746294; 10.30509461847965; 10.30488294090483; 10.30422590242215; 10.30384905150036; 10.30208292513454; 10.30135162625518; 10.30106083753661; 10.30098471931155; 10.30141604606089; 10.30441355713189; 10.30751972107803; 10.3057573960063; 10.30138651795233; 10.30164011417437) # column D
    5 = @(15.37566571272588; 15.41272278575937; 15.43712527500929; 15.44748492608324; 15.4492302505881; 15.43726330583842; 15.38810110731427; 15.30474879138554; 15.25142452617461; 15.22887527038297; 15.22058138217472; 15.22235673157878; 15.22818802004026; 15.23179174397227; 15.25293249424779; 15.2663387379266; 15.27421047704362; 15.27690335279992; 15.26489498035912; 15.22646858422748; 15.20278262625482; 15.21529381595976; 15.26554719130968; 15.32590475052233) # column E
    6 = @(30.63715244501142; 30.70394629571991; 30.75098293579243; 30.77166455135157; 30.77519010738467; 30.75125572774859; 30.6589318250496; 30.52573649220693; 30.45710036338796; 30.43222670724759; 30.4237206063087; 30.42551194435557; 30.43150860823662; 30.43530063361456; 30.4588536857793; 30.47492901848297; 30.48477276777222; 30.48820833276281; 30.47315591706196; 30.42972246542119; 30.40665769966514; 30.41848096621918; 30.4739556615916; 30.55664013854091) # column F
    7 = @(29.982748526454; 30.038353843926; 30.08040737683832; 30.09952993064621; 30.10282499694162; 30.08065723736512; 30.00027695383073; 29.90557817105891; 29.87452712062746; 29.86878438902997; 29.86781515013134; 29.86797029164228; 29.86868049412455; 29.86927247836751; 29.87507112119905; 29.88077566316392; 29.88484582582613; 29.88635941999203; 29.88008673716263; 29.86843918038676; 29.86785223916746; 29.86752290585139; 29.8803957381467; 29.92443999434404) # column G
    8 = @(14.68347457919492; 14.72314415169157; 14.74946830123815; 14.76069052071921; 14.76258386297943; 14.74961764368765; 14.69674462312881; 14.60865227164001; 14.55341546850387; 14.53034197820495; 14.52189968683589; 14.52370476260139; 14.52964151210142; 14.53331637159028; 14.55496468810718; 14.56877115108239; 14.57690561015218; 14.57969301250638; 14.5672814204836; 14.5278897367049; 14.50386518019956; 14.51653021666281; 14.56795431421836; 14.63081650360562) # column H
    9 = @(25.74724792094496; 25.84644628663432; 25.91101729694153; 25.9382528786039; 25.94283107451473; 25.91138086909245; 25.78069209867862; 25.55341401184412; 25.40402851155162; 25.33987232149124; 25.31612325306756; 25.32121379392919; 25.33790754426765; 25.34820395940777; 25.40829765366439; 25.44613577938252; 25.46825699949246; 25.47580833465285; 25.44207082947961; 25.33298938814746; 25.26487746610644; 25.30093949800455; 25.44390744828999; 25.61180275610493) # column I
    10 = @(10.95712308390268; 10.97771670818454; 10.99105481800493; 10.99666507286497; 10.99760722678456; 10.99112977120583; 10.96408010463011; 10.91651723906108; 10.88488431255016; 10.87120634234678; 10.86612875676236; 10.86721777741102; 10.87078656519235; 10.87298581717346; 10.88579249139954; 10.89383103250347; 10.89852163108581; 10.90012131464611; 10.89296838051231; 10.86973556174831; 10.85514570862973; 10.86287836221811; 10.8933581701033; 10.92880055704584) # column J
    11 = @(8.511301873186509; 8.247175234760512; 8.081300680277812; 8.012873530798641; 8.001464035739506; 8.080381080533549; 8.421055526269207; 9.12743310422313; 9.656876906412741; 9.886766189016686; 9.972207353002061; 9.953878336462232; 9.893828003106012; 9.856834384006143; 9.641629290624538; 9.506773078618103; 9.428177818583553; 9.401391223541546; 9.521235594018467; 9.911510283997059; 10.15716270027624; 10.02692566650586; 9.5147004040195; 8.923027299042349) # column K
    12 = @(9.474218825952276; 9.435719449368428; 9.413146093767612; 9.404221429073568; 9.402756231124373; 9.413024614499882; 9.460726946311654; 9.56243874991811; 9.641747295686738; 9.678731908237216; 9.692859507525226; 9.689811565769627; 9.679891789449655; 9.673831328739597; 9.639347825416611; 9.618419875897786; 9.606468515354253; 9.6024369763072; 9.620638865320394; 9.682802213318748; 9.724138729409965; 9.702014469007688; 9.61963540916711; 9.534088757581411) # column L
    14 = @(19.45540016240225; 19.50525215757478; 19.53744588314023; 19.55096449475021; 19.55323340254509; 19.53762658135021; 19.47226088957532; 19.35660462607769; 19.27920331115939; 19.24562182608338; 19.23313860425893; 19.23581672427608; 19.24459015193047; 19.24999449458275; 19.28143064153943; 19.30113225095873; 19.31261745267226; 19.31653251015944; 19.29901911536479; 19.24200685705233; 19.20610585020881; 19.22514275088326; 19.29997397023571; 19.38655848229841) # column N
    15 = @(22.49773452961068; 22.56011322036555; 22.6024345154389; 22.62069108283584; 22.62378356968821; 22.60267664020959; 22.51840784163933; 22.38508741690259; 22.30663979584242; 22.2751933620606; 22.26389550760284; 22.26630155488115; 22.27425164807657; 22.27920079414025; 22.30878025881352; 22.32801263438702; 22.33947362928282; 22.3434226426; 22.32592401161964; 22.27189994765723; 22.24014916116861; 22.25676952401314; 22.32686701944653; 22.4177315990006) # column O
}

foreach ($col in $data.Keys) {
    $colValues = $data[$col]
    for ($i = 0; $i -lt $colValues.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, $col).Value2 = $colValues[$i]
    }
}
